$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated error-table values (naive component forecaster bug fix)
# Row 2 (Q0, horizon 6)
$ws.Range("B2").Value = -0.011922918044421
$ws.Range("C2").Value = 1.502158260014067
$ws.Range("D2").Value = 8.459308576868294
$ws.Range("E2").Value = 2.908489053936475
$ws.Range("F2").Value = 2.973831532800685
$ws.Range("G2").Value = 23

# Row 3 (Q1, horizon 7)
$ws.Range("B3").Value = 0.5658384863776997
$ws.Range("C3").Value = 2.129283520569611
$ws.Range("D3").Value = 18.65469822571722
$ws.Range("E3").Value = 4.319108498951747
$ws.Range("F3").Value = 4.382647333937844
$ws.Range("G3").Value = 22

# Row 4 (Q2, horizon 8)
$ws.Range("B4").Value = -0.9240374313429358
$ws.Range("C4").Value = 1.444924223990262
$ws.Range("D4").Value = 8.174846682920178
$ws.Range("E4").Value = 2.85916887974813
$ws.Range("F4").Value = 2.772553260771954
$ws.Range("G4").Value = 21

# Row 5 (Q3, horizon 9)
$ws.Range("B5").Value = 0.07040277154089039
$ws.Range("C5").Value = 0.5120019363376737
$ws.Range("D5").Value = 0.4676573765475683
$ws.Range("E5").Value = 0.6838547920045368
$ws.Range("F5").Value = 0.6978921846357898
$ws.Range("G5").Value = 20

# Row 6 (Q4, horizon 10)
$ws.Range("B6").Value = 0.1483703450672545
$ws.Range("C6").Value = 0.8649345531276215
$ws.Range("D6").Value = 1.90040254551332
$ws.Range("E6").Value = 1.378550886080496
$ws.Range("F6").Value = 1.408099328451515
$ws.Range("G6").Value = 19

# Row 7 (Q5, horizon 11)
$ws.Range("B7").Value = -0.07978550902713241
$ws.Range("C7").Value = 0.7076988813810415
$ws.Range("D7").Value = 1.325148457270131
$ws.Range("E7").Value = 1.151150927233319
$ws.Range("F7").Value = 1.181676006574844
$ws.Range("G7").Value = 18

# Row 8 (Q6, horizon 12)
$ws.Range("B8").Value = -0.004526189646546301
$ws.Range("C8").Value = 0.5864325954492303
$ws.Range("D8").Value = 0.5709791096934894
$ws.Range("E8").Value = 0.7556315965425807
$ws.Range("F8").Value = 0.7788732485180572
$ws.Range("G8").Value = 17

# Row 9 (Q7, horizon 13)
$ws.Range("B9").Value = 0.2502318746956579
$ws.Range("C9").Value = 0.5174934839310298
$ws.Range("D9").Value = 0.4275779196144051
$ws.Range("E9").Value = 0.653894425434569
$ws.Range("F9").Value = 0.6239332686546548
$ws.Range("G9").Value = 16

# Row 10 (Q8, horizon 14)
$ws.Range("B10").Value = 0.1901465307434204
$ws.Range("C10").Value = 0.4556814281592019
$ws.Range("D10").Value = 0.3317531417364475
$ws.Range("E10").Value = 0.575980157415555
$ws.Range("F10").Value = 0.5627713046510059
$ws.Range("G10").Value = 15

# Row 11 (Q9, horizon 15) - G11 unchanged (stays 13)
$ws.Range("B11").Value = 0.2349662137414842
$ws.Range("C11").Value = 0.4912372055275555
$ws.Range("D11").Value = 0.3821517403275286
$ws.Range("E11").Value = 0.6181842284687702
$ws.Range("F11").Value = 0.5951368220458327
